$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 157-158: second half of the masked-word placeholders (words 8 and 9 of
# sentence 7) gain the same red "masked" style as the earlier placeholder rows
# (149-156), including blank, styled E/F cells.
$ws.Range("D157:F157").Interior.Color = 255
$ws.Range("D158:F158").Interior.Color = 255

# Rows 159-168: reveal words 10-19 of sentence 7 ("investigation and submits
# that report to the police accountability board"), each with a shared-string
# word in column E and the constant count 1 in column F.
$words = @(
    @{ Row = 159; Word = "[b'investigation']" },
    @{ Row = 160; Word = "[b'and']" },
    @{ Row = 161; Word = "[b'submits']" },
    @{ Row = 162; Word = "[b'that']" },
    @{ Row = 163; Word = "[b'report']" },
    @{ Row = 164; Word = "[b'to']" },
    @{ Row = 165; Word = "[b'the']" },
    @{ Row = 166; Word = "[b'police']" },
    @{ Row = 167; Word = "[b'accountability']" },
    @{ Row = 168; Word = "[b'board']" }
)

foreach ($entry in $words) {
    $r = $entry.Row
    $ws.Cells.Item($r, 5).Value = $entry.Word
    $ws.Cells.Item($r, 6).Value = 1
}

# Update the view state: scrolled up a few rows and the active selection now
# sits on the very last cell that was just filled in.
$ws.Range("F168").Select()
$excel.ActiveWindow.ScrollRow = 139
$excel.ActiveWindow.ScrollColumn = 3
